$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("AA5").Value = 2
$ws.Range("AA27").Value = 712
$ws.Range("AA28").Value = 712
$ws.Range("AA50").Value = 3020
$ws.Range("AB50").Value = 25774
$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("AA21").Value = 1
$ws.Range("AB21").Value = 18
$ws.Range("AA23").Value = 1
$ws.Range("AA24").Value = 2
$ws.Range("AB24").Value = 40
$ws.Range("AA26").Value = 4
$ws.Range("AB26").Value = 69
$ws.Range("AA27").Value = 5
$ws.Range("AB27").Value = 86
$ws.Range("AA28").Value = 6
$ws.Range("AB28").Value = 105
$ws.Range("AA29").Value = 6
$ws.Range("AB29").Value = 137
$ws.Range("AA30").Value = 8
$ws.Range("AB30").Value = 158
$ws.Range("AA31").Value = 10
$ws.Range("AB31").Value = 190
$ws.Range("AA32").Value = 13
$ws.Range("AB32").Value = 232
$ws.Range("AA33").Value = 15
$ws.Range("AB33").Value = 266
$ws.Range("AA34").Value = 19
$ws.Range("AB34").Value = 306
$ws.Range("AA35").Value = 22
$ws.Range("AB35").Value = 342
$ws.Range("AA36").Value = 24
$ws.Range("AB36").Value = 397
$ws.Range("AA37").Value = 26
$ws.Range("AB37").Value = 467
$ws.Range("AA38").Value = 35
$ws.Range("AB38").Value = 524
$ws.Range("AA39").Value = 38
$ws.Range("AB39").Value = 574
$ws.Range("AA40").Value = 40
$ws.Range("AB40").Value = 633
$ws.Range("AA41").Value = 48
$ws.Range("AB41").Value = 704
$ws.Range("AA42").Value = 52
$ws.Range("AB42").Value = 753
$ws.Range("AA43").Value = 54
$ws.Range("AB43").Value = 808
$ws.Range("AA44").Value = 56
$ws.Range("AB44").Value = 863
$ws.Range("AA45").Value = 62
$ws.Range("AB45").Value = 931
$ws.Range("AA46").Value = 63
$ws.Range("AB46").Value = 983
$ws.Range("AA47").Value = 67
$ws.Range("AB47").Value = 1040
$ws.Range("AA48").Value = 72
$ws.Range("AB48").Value = 1080
$ws.Range("AA49").Value = 74
$ws.Range("AB49").Value = 1126
$ws.Range("AA50").Value = 77
$ws.Range("AB50").Value = 1156
$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("AA25").Value = 60
$ws.Range("AB25").Value = 631
$ws.Range("AA26").Value = 78
$ws.Range("AB26").Value = 736
$ws.Range("AA27").Value = 83
$ws.Range("AB27").Value = 834
$ws.Range("AA28").Value = 110
$ws.Range("AB28").Value = 1031
$ws.Range("AA29").Value = 127
$ws.Range("AB29").Value = 1139
$ws.Range("AA30").Value = 138
$ws.Range("AB30").Value = 1286
$ws.Range("AA31").Value = 146
$ws.Range("AB31").Value = 1385
$ws.Range("AA32").Value = 156
$ws.Range("AB32").Value = 1557
$ws.Range("AA33").Value = 171
$ws.Range("AB33").Value = 1740
$ws.Range("AA36").Value = 197
$ws.Range("AB36").Value = 2126
$ws.Range("AA37").Value = 196
$ws.Range("AB37").Value = 2165
$ws = $wb.Worksheets.Item("Ventilated")
$ws.Range("AA25").Value = 12
$ws.Range("AB25").Value = 35
$ws.Range("AA26").Value = 15
$ws.Range("AB26").Value = 43
$ws.Range("AA27").Value = 23
$ws.Range("AB27").Value = 91
$ws.Range("AA28").Value = 22
$ws.Range("AB28").Value = 109
$ws.Range("AA29").Value = 27
$ws.Range("AB29").Value = 126
$ws.Range("AA30").Value = 26
$ws.Range("AB30").Value = 132
$ws.Range("AA31").Value = 28
$ws.Range("AB31").Value = 152
$ws.Range("AA32").Value = 35
$ws.Range("AB32").Value = 164
$ws.Range("AA33").Value = 38
$ws.Range("AB33").Value = 174
$ws.Range("AA36").Value = 46
$ws.Range("AB36").Value = 239
$ws.Range("AA37").Value = 50
$ws.Range("AB37").Value = 251
$ws = $wb.Worksheets.Item("Tested")
$ws.Range("AA9").ClearContents()
$ws.Range("AB9").Value = 1244
$ws.Range("AB10").Value = 1353
$ws.Range("AB11").Value = 1417
$ws.Range("AB12").Value = 1517
$ws.Range("AB13").Value = 1594
$ws.Range("AB14").Value = 1667
$ws.Range("AA15").ClearContents()
$ws.Range("AB15").Value = 1780
$ws.Range("AB16").Value = 1992
$ws.Range("AB17").Value = 2238
$ws.Range("AB18").Value = 2603
$ws.Range("AB19").Value = 3081
$ws.Range("AB20").Value = 3731
$ws.Range("AB21").Value = 4175
$ws.Range("AB22").Value = 4886
$ws.Range("AB23").Value = 5533
$ws.Range("AB24").Value = 6273
$ws.Range("AB25").Value = 6987
$ws.Range("AB26").Value = 7513
$ws.Range("AB27").Value = 7878
$ws.Range("AB28").Value = 8135
$ws.Range("AB29").Value = 8537
$ws.Range("AB30").Value = 8874
$ws.Range("AB31").Value = 9313
$ws.Range("AB32").Value = 9845
$ws.Range("AB33").Value = 10346
$ws.Range("AB34").Value = 10949
$ws.Range("AB35").Value = 11271
$ws.Range("AB36").Value = 11949
$ws.Range("AB37").Value = 12703
$ws.Range("AB38").Value = 13356
$ws.Range("AB39").Value = 14289
$ws.Range("AB40").Value = 15106
$ws.Range("AB41").Value = 15719
$ws.Range("AB42").Value = 16049
$ws.Range("AB43").Value = 16643
$ws.Range("AB44").Value = 17275
$ws.Range("AB45").Value = 17982
$ws.Range("AB46").Value = 18470
$ws.Range("AB47").Value = 18916
$ws.Range("AB48").Value = 19261
$ws.Range("AB49").Value = 19431
$ws.Range("AB50").Value = 19702
$ws.Range("AB51").Value = 19702
$ws.Range("AB52").Value = 19702
